# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (F column) and "最低票价" (G column) figures
# scraped for each event, across all four sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 151
$ws1.Cells.Item(2, 7).Value = 73
$ws1.Cells.Item(3, 6).Value = 1345
$ws1.Cells.Item(3, 7).Value = 73
$ws1.Cells.Item(4, 6).Value = 1157
$ws1.Cells.Item(5, 6).Value = 1046
$ws1.Cells.Item(6, 6).Value = 1838
$ws1.Cells.Item(7, 6).Value = 584
$ws1.Cells.Item(8, 6).Value = 1221
$ws1.Cells.Item(11, 6).Value = 130
$ws1.Cells.Item(13, 6).Value = 89
$ws1.Cells.Item(15, 6).Value = 732
$ws1.Cells.Item(16, 6).Value = 192
$ws1.Cells.Item(21, 6).Value = 175
$ws1.Cells.Item(23, 6).Value = 52
$ws1.Cells.Item(25, 6).Value = 171
$ws1.Cells.Item(27, 6).Value = 887
$ws1.Cells.Item(28, 6).Value = 328
$ws1.Cells.Item(29, 6).Value = 170
$ws1.Cells.Item(30, 6).Value = 53

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(7, 6).Value = 263
$ws2.Cells.Item(11, 6).Value = 123
$ws2.Cells.Item(12, 6).Value = 24

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 317

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 317
$ws4.Cells.Item(3, 6).Value = 151
$ws4.Cells.Item(3, 7).Value = 73
$ws4.Cells.Item(4, 6).Value = 1345
$ws4.Cells.Item(4, 7).Value = 73
$ws4.Cells.Item(5, 6).Value = 1157
$ws4.Cells.Item(7, 6).Value = 1838
$ws4.Cells.Item(8, 6).Value = 584
$ws4.Cells.Item(9, 6).Value = 1221
$ws4.Cells.Item(13, 6).Value = 130
$ws4.Cells.Item(15, 6).Value = 89
$ws4.Cells.Item(17, 6).Value = 732
$ws4.Cells.Item(18, 6).Value = 192
$ws4.Cells.Item(27, 6).Value = 263
$ws4.Cells.Item(28, 6).Value = 263
$ws4.Cells.Item(29, 6).Value = 175
$ws4.Cells.Item(31, 6).Value = 52
$ws4.Cells.Item(33, 6).Value = 171
$ws4.Cells.Item(35, 6).Value = 887
$ws4.Cells.Item(36, 6).Value = 328
$ws4.Cells.Item(39, 6).Value = 170
$ws4.Cells.Item(40, 6).Value = 53
$ws4.Cells.Item(43, 6).Value = 123
$ws4.Cells.Item(44, 6).Value = 123
$ws4.Cells.Item(47, 6).Value = 24

Write-Host "edit.ps1 applied"
